$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, pushing existing rows 19-46 down to 20-47.
# Excel's default Insert() copies formatting from the row above, which
# matches the date-style (s="2") needed on the new D19 cell.
$ws.Rows.Item(19).Insert()

# Fill in the carry-over (unchanged) descriptive columns for the new record,
# copied from the same values as the rest of this market/product block.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value2 = 44589
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103002
$ws.Range("J19").Value = "Ciruela"
$ws.Range("K19").Value = "Black Amber"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 220
$ws.Range("N19").Value = 7500
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 7773
$ws.Range("Q19").Value = "$/caja 16 kilos granel"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 486
$ws.Range("T19").Value = 16
